# Fills in self-assessment scores that were left blank for a few
# group members across the four evaluation sheets.

$wb = $excel.ActiveWorkbook

# --- "Group and Self Assessment": student 1212047 (row 12) scores ---
$ws1 = $wb.Worksheets.Item("Group and Self Assessment")
$ws1.Range("D12").Value = 4
$ws1.Range("E12").Value = 3
$ws1.Range("F12").Value = 4
$ws1.Range("G12").Value = 5

# --- "User Stories": student 1221720 (row 6) score ---
$ws2 = $wb.Worksheets.Item("User Stories")
$ws2.Range("C6").Value = 4

# --- "Project Development": self-assessment column E for rows 4-7 ---
$ws3 = $wb.Worksheets.Item("Project Development")
$ws3.Range("E4").Value = 5
$ws3.Range("E5").Value = 4
$ws3.Range("E6").Value = 3
$ws3.Range("E7").Value = 4

# --- "Project Management": self-assessment column E for rows 4-14 ---
$ws4 = $wb.Worksheets.Item("Project Management")
$ws4.Range("E4").Value = 5
$ws4.Range("E5").Value = 4
$ws4.Range("E6").Value = 5
$ws4.Range("E7").Value = 5
$ws4.Range("E8").Value = 5
$ws4.Range("E9").Value = 4
$ws4.Range("E10").Value = 5
$ws4.Range("E11").Value = 3
$ws4.Range("E12").Value = 4
$ws4.Range("E13").Value = 4
$ws4.Range("E14").Value = 4
